$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D; this shifts the existing D:K data to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting from column F (the shifted original column D) into the new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write the refreshed/restated quarterly figures across columns D through M
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 43465
$arr[0,1] = 43373
$arr[0,2] = 43281
$arr[0,3] = 43190
$arr[0,4] = 43100
$arr[0,5] = 43008
$arr[0,6] = 42916
$arr[0,7] = 42825
$arr[0,8] = 42735
$arr[0,9] = 42643
$ws.Range("D7:M7").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 3803800
$arr[0,1] = 4191400
$arr[0,2] = 3534500
$arr[0,3] = 3340300
$arr[0,4] = 3237700
$arr[0,5] = 3118100
$arr[0,6] = 3043300
$arr[0,7] = 3125700
$arr[0,8] = 2988400
$arr[0,9] = 2972300
$ws.Range("D8:M8").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 3297700
$arr[0,1] = 3259700
$arr[0,2] = 2692900
$arr[0,3] = 2647700
$arr[0,4] = 2405700
$arr[0,5] = 2349100
$arr[0,6] = 2319300
$arr[0,7] = 2210800
$arr[0,8] = 2269000
$arr[0,9] = 2174500
$ws.Range("D9:M9").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 506100
$arr[0,1] = 931700
$arr[0,2] = 841600
$arr[0,3] = 692700
$arr[0,4] = 831900
$arr[0,5] = 769000
$arr[0,6] = 724100
$arr[0,7] = 914900
$arr[0,8] = 719400
$arr[0,9] = 797900
$ws.Range("D10:M10").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 15700
$arr[0,1] = 13600
$arr[0,2] = 12000
$arr[0,3] = 10000
$arr[0,4] = 14500
$arr[0,5] = 10000
$arr[0,6] = 9900
$arr[0,7] = 8400
$arr[0,8] = 10500
$arr[0,9] = 9200
$ws.Range("D12:M12").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D13:M13").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -25900
$arr[0,1] = -20400
$arr[0,2] = -11500
$arr[0,3] = -7100
$arr[0,4] = 15100
$arr[0,5] = 10600
$arr[0,6] = -65200
$arr[0,7] = -700
$arr[0,8] = 26900
$arr[0,9] = "NA"
$ws.Range("D14:M14").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D15:M15").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 3505700
$arr[0,1] = 3476000
$arr[0,2] = 2909800
$arr[0,3] = 2849600
$arr[0,4] = 2671600
$arr[0,5] = 2622900
$arr[0,6] = 2450700
$arr[0,7] = 2404300
$arr[0,8] = 3295900
$arr[0,9] = 2402600
$ws.Range("D17:M17").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 298200
$arr[0,1] = 715400
$arr[0,2] = 624700
$arr[0,3] = 490700
$arr[0,4] = 566000
$arr[0,5] = 495200
$arr[0,6] = 592600
$arr[0,7] = 721500
$arr[0,8] = -307500
$arr[0,9] = 569700
$ws.Range("D18:M18").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -147900
$arr[0,1] = -97500
$arr[0,2] = -408900
$arr[0,3] = -3700
$arr[0,4] = -365400
$arr[0,5] = -99500
$arr[0,6] = -26500
$arr[0,7] = 46800
$arr[0,8] = -229800
$arr[0,9] = -283500
$ws.Range("D20:M20").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 340500
$arr[0,1] = 819900
$arr[0,2] = 400600
$arr[0,3] = 676700
$arr[0,4] = 395300
$arr[0,5] = 586700
$arr[0,6] = 751200
$arr[0,7] = 942400
$arr[0,8] = -356700
$arr[0,9] = 462600
$ws.Range("D21:M21").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 132000
$arr[0,1] = 141300
$arr[0,2] = 140100
$arr[0,3] = 121100
$arr[0,4] = 131700
$arr[0,5] = 141600
$arr[0,6] = 147200
$arr[0,7] = 142400
$arr[0,8] = 159500
$arr[0,9] = "NA"
$ws.Range("D22:M22").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 18300
$arr[0,1] = 476600
$arr[0,2] = 75600
$arr[0,3] = 365900
$arr[0,4] = 68800
$arr[0,5] = 254200
$arr[0,6] = 418900
$arr[0,7] = 625800
$arr[0,8] = -696900
$arr[0,9] = 286200
$ws.Range("D23:M23").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 64100
$arr[0,1] = 106900
$arr[0,2] = -50700
$arr[0,3] = 70800
$arr[0,4] = -11400
$arr[0,5] = 58400
$arr[0,6] = 126000
$arr[0,7] = 153200
$arr[0,8] = -46600
$arr[0,9] = 84600
$ws.Range("D24:M24").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D25:M25").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -45800
$arr[0,1] = 369700
$arr[0,2] = 126300
$arr[0,3] = 295100
$arr[0,4] = 80200
$arr[0,5] = 195800
$arr[0,6] = 292900
$arr[0,7] = 472700
$arr[0,8] = -650200
$arr[0,9] = 201700
$ws.Range("D26:M26").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -20100
$arr[0,1] = 344600
$arr[0,2] = 140300
$arr[0,3] = 270100
$arr[0,4] = 99000
$arr[0,5] = 204900
$arr[0,6] = 279400
$arr[0,7] = 446200
$arr[0,8] = -628800
$arr[0,9] = 219400
$ws.Range("D27:M27").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D28:M28").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 2200
$arr[0,8] = 900
$arr[0,9] = 1300
$ws.Range("D29:M29").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D30:M30").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D31:M31").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 147900
$arr[0,1] = 97500
$arr[0,2] = 408900
$arr[0,3] = 3700
$arr[0,4] = 365400
$arr[0,5] = 99500
$arr[0,6] = 26500
$arr[0,7] = -46800
$arr[0,8] = 229800
$arr[0,9] = 283500
$ws.Range("D32:M32").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -20100
$arr[0,1] = 344600
$arr[0,2] = 140300
$arr[0,3] = 270100
$arr[0,4] = 99000
$arr[0,5] = 204900
$arr[0,6] = 279400
$arr[0,7] = 448400
$arr[0,8] = -628000
$arr[0,9] = 220700
$ws.Range("D33:M33").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D34:M34").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -20100
$arr[0,1] = 344600
$arr[0,2] = 140300
$arr[0,3] = 270100
$arr[0,4] = 99000
$arr[0,5] = 204900
$arr[0,6] = 279400
$arr[0,7] = 448400
$arr[0,8] = -628000
$arr[0,9] = 220700
$ws.Range("D35:M35").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 43465
$arr[0,1] = 43373
$arr[0,2] = 43281
$arr[0,3] = 43190
$arr[0,4] = 43100
$arr[0,5] = 43008
$arr[0,6] = 42916
$arr[0,7] = 42825
$arr[0,8] = 42735
$arr[0,9] = 42643
$ws.Range("D38:M38").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1422300
$arr[0,1] = 1254000
$arr[0,2] = 1320300
$arr[0,3] = 875100
$arr[0,4] = 967900
$arr[0,5] = 1397900
$arr[0,6] = 1464300
$arr[0,7] = 1641600
$arr[0,8] = 1662600
$arr[0,9] = 2034300
$ws.Range("D41:M41").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 604400
$arr[0,1] = 533500
$arr[0,2] = 527100
$arr[0,3] = 625600
$arr[0,4] = 590400
$arr[0,5] = 585700
$arr[0,6] = 686500
$arr[0,7] = 250900
$arr[0,8] = 295300
$arr[0,9] = 107400
$ws.Range("D42:M42").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1211400
$arr[0,1] = 1512100
$arr[0,2] = 1158800
$arr[0,3] = 1130500
$arr[0,4] = 1263800
$arr[0,5] = 1077900
$arr[0,6] = 1019500
$arr[0,7] = 752400
$arr[0,8] = 658900
$arr[0,9] = 796900
$ws.Range("D43:M43").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 2175800
$arr[0,1] = 2293300
$arr[0,2] = 2080300
$arr[0,3] = 1712500
$arr[0,4] = 1755400
$arr[0,5] = 1473200
$arr[0,6] = 1514600
$arr[0,7] = 1375900
$arr[0,8] = 1299400
$arr[0,9] = 1209100
$ws.Range("D44:M44").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 68500
$arr[0,1] = 34200
$arr[0,2] = 40300
$arr[0,3] = 34300
$arr[0,4] = 35400
$arr[0,5] = 51200
$arr[0,6] = 52400
$arr[0,7] = 111400
$arr[0,8] = 116600
$arr[0,9] = 39000
$ws.Range("D45:M45").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 5482400
$arr[0,1] = 5627100
$arr[0,2] = 5126800
$arr[0,3] = 4378000
$arr[0,4] = 4612900
$arr[0,5] = 4585900
$arr[0,6] = 4737400
$arr[0,7] = 4132100
$arr[0,8] = 4032900
$arr[0,9] = 4186700
$ws.Range("D46:M46").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 485700
$arr[0,1] = 330300
$arr[0,2] = 347500
$arr[0,3] = 339400
$arr[0,4] = 329600
$arr[0,5] = 378500
$arr[0,6] = 380000
$arr[0,7] = 362800
$arr[0,8] = 345300
$arr[0,9] = 426900
$ws.Range("D47:M47").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 8142600
$arr[0,1] = 8395600
$arr[0,2] = 8108300
$arr[0,3] = 7782700
$arr[0,4] = 7630300
$arr[0,5] = 7640900
$arr[0,6] = 7806300
$arr[0,7] = 7344800
$arr[0,8] = 7277900
$arr[0,9] = 7403700
$ws.Range("D48:M48").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 702700
$arr[0,1] = 704600
$arr[0,2] = 702400
$arr[0,3] = 695500
$arr[0,4] = 699300
$arr[0,5] = 704400
$arr[0,6] = 714300
$arr[0,7] = 689900
$arr[0,8] = 696900
$arr[0,9] = 699100
$ws.Range("D49:M49").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D50:M50").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D51:M51").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 362700
$arr[0,1] = 528000
$arr[0,2] = 625000
$arr[0,3] = 343700
$arr[0,4] = 403600
$arr[0,5] = 326100
$arr[0,6] = 353300
$arr[0,7] = 386600
$arr[0,8] = 503100
$arr[0,9] = 457300
$ws.Range("D52:M52").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D53:M53").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 15176200
$arr[0,1] = 15585600
$arr[0,2] = 14910100
$arr[0,3] = 13539200
$arr[0,4] = 13675700
$arr[0,5] = 13635900
$arr[0,6] = 13991200
$arr[0,7] = 12916200
$arr[0,8] = 12856000
$arr[0,9] = 13173700
$ws.Range("D54:M54").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 2138500
$arr[0,1] = 2405700
$arr[0,2] = 2117200
$arr[0,3] = 1453100
$arr[0,4] = 1350000
$arr[0,5] = 1233400
$arr[0,6] = 1348000
$arr[0,7] = 1257900
$arr[0,8] = 1623700
$arr[0,9] = 1777600
$ws.Range("D57:M57").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 2889300
$arr[0,1] = 3098000
$arr[0,2] = 3108800
$arr[0,3] = 2612000
$arr[0,4] = 2795400
$arr[0,5] = 3374100
$arr[0,6] = 3695000
$arr[0,7] = 3205900
$arr[0,8] = 3233000
$arr[0,9] = 731800
$ws.Range("D58:M58").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 898600
$arr[0,1] = 723300
$arr[0,2] = 691000
$arr[0,3] = 577600
$arr[0,4] = 761100
$arr[0,5] = 719900
$arr[0,6] = 760100
$arr[0,7] = 881200
$arr[0,8] = 858600
$arr[0,9] = 814100
$ws.Range("D59:M59").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 5926500
$arr[0,1] = 6227000
$arr[0,2] = 5917000
$arr[0,3] = 4642700
$arr[0,4] = 4906500
$arr[0,5] = 5327400
$arr[0,6] = 5803000
$arr[0,7] = 5344900
$arr[0,8] = 5715300
$arr[0,9] = 3323500
$ws.Range("D60:M60").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 6194300
$arr[0,1] = 6150000
$arr[0,2] = 6230900
$arr[0,3] = 5469200
$arr[0,4] = 5685600
$arr[0,5] = 4740400
$arr[0,6] = 4970200
$arr[0,7] = 4871100
$arr[0,8] = 5144300
$arr[0,9] = 7597500
$ws.Range("D61:M61").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1539800
$arr[0,1] = 1564800
$arr[0,2] = 1519100
$arr[0,3] = 1546400
$arr[0,4] = 1624600
$arr[0,5] = 1772800
$arr[0,6] = 1746400
$arr[0,7] = 1575400
$arr[0,8] = 1569400
$arr[0,9] = 1162700
$ws.Range("D62:M62").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D63:M63").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D64:M64").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D65:M65").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 13435900
$arr[0,1] = 13746600
$arr[0,2] = 13436000
$arr[0,3] = 11491200
$arr[0,4] = 12004700
$arr[0,5] = 11654000
$arr[0,6] = 12338600
$arr[0,7] = 11595600
$arr[0,8] = 12176600
$arr[0,9] = 11845300
$ws.Range("D66:M66").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D68:M68").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D69:M69").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D70:M70").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D71:M71").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1257700
$arr[0,1] = 1447200
$arr[0,2] = 1100700
$arr[0,3] = 1343100
$arr[0,4] = 1071200
$arr[0,5] = 1226800
$arr[0,6] = 1020100
$arr[0,7] = 714900
$arr[0,8] = 264700
$arr[0,9] = 914200
$ws.Range("D72:M72").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D73:M73").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D74:M74").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D75:M75").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1740200
$arr[0,1] = 1839100
$arr[0,2] = 1474100
$arr[0,3] = 2048000
$arr[0,4] = 1671000
$arr[0,5] = 1981900
$arr[0,6] = 1652600
$arr[0,7] = 1320600
$arr[0,8] = 679400
$arr[0,9] = 1328400
$ws.Range("D76:M76").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D77:M77").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 43465
$arr[0,1] = 43373
$arr[0,2] = 43281
$arr[0,3] = 43190
$arr[0,4] = 43100
$arr[0,5] = 43008
$arr[0,6] = 42916
$arr[0,7] = 42825
$arr[0,8] = 42735
$arr[0,9] = 42643
$ws.Range("D80:M80").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -20100
$arr[0,1] = 344600
$arr[0,2] = 140300
$arr[0,3] = 270100
$arr[0,4] = 99000
$arr[0,5] = 204900
$arr[0,6] = 279400
$arr[0,7] = 448400
$arr[0,8] = -628000
$arr[0,9] = 220700
$ws.Range("D81:M81").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 190100
$arr[0,1] = 202100
$arr[0,2] = 184800
$arr[0,3] = 189700
$arr[0,4] = 194800
$arr[0,5] = 191000
$arr[0,6] = 185100
$arr[0,7] = 174200
$arr[0,8] = 177700
$arr[0,9] = 176400
$ws.Range("D83:M83").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D84:M84").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D85:M85").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D86:M86").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D87:M87").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D88:M88").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 258400
$arr[0,1] = 542100
$arr[0,2] = 1109400
$arr[0,3] = 461700
$arr[0,4] = 341000
$arr[0,5] = 377300
$arr[0,6] = -232900
$arr[0,7] = 141000
$arr[0,8] = 179800
$arr[0,9] = 704800
$ws.Range("D89:M89").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -231000
$arr[0,1] = -169200
$arr[0,2] = -189900
$arr[0,3] = -103700
$arr[0,4] = -187200
$arr[0,5] = -198300
$arr[0,6] = -120300
$arr[0,7] = -67700
$arr[0,8] = -418100
$arr[0,9] = -262200
$ws.Range("D91:M91").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D92:M92").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D93:M93").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -216700
$arr[0,1] = -148700
$arr[0,2] = -169100
$arr[0,3] = -103400
$arr[0,4] = -341100
$arr[0,5] = -195500
$arr[0,6] = -9900
$arr[0,7] = -68200
$arr[0,8] = -206800
$arr[0,9] = -130100
$ws.Range("D94:M94").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = -384500
$arr[0,3] = 0
$arr[0,4] = -256100
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = -247800
$arr[0,9] = 0
$ws.Range("D96:M96").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D97:M97").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D98:M98").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D99:M99").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 105600
$arr[0,1] = -412800
$arr[0,2] = -428300
$arr[0,3] = -444800
$arr[0,4] = -416900
$arr[0,5] = -275700
$arr[0,6] = 35200
$arr[0,7] = -105200
$arr[0,8] = -237800
$arr[0,9] = -197800
$ws.Range("D100:M100").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 20900
$arr[0,1] = -46900
$arr[0,2] = -66800
$arr[0,3] = -6200
$arr[0,4] = -13000
$arr[0,5] = 27500
$arr[0,6] = -24500
$arr[0,7] = 11300
$arr[0,8] = 14900
$arr[0,9] = -10500
$ws.Range("D101:M101").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 168300
$arr[0,1] = -66300
$arr[0,2] = 445200
$arr[0,3] = -92800
$arr[0,4] = -430000
$arr[0,5] = -66400
$arr[0,6] = -232200
$arr[0,7] = -21000
$arr[0,8] = -254500
$arr[0,9] = 362100
$ws.Range("D102:M102").Value = $arr

